$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "27.092.81"
$ws.Range("E2").Value = "  +0.51%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.680.39"
$ws.Range("E3").Value = "  +0.49%  "
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "215.13"
$ws.Range("E5").Value = "  +0.07%  "
$ws.Range("E6").Value = "  -0.03%  "
$ws.Range("E7").Value = "  +0.03%  "
$ws.Range("E8").Value = "  +2.12%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "21.31"
$ws.Range("E9").Value = "  +5.68%  "
$ws.Range("E10").Value = "  +0.47%  "
$ws.Range("E11").Value = "  -0.49%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.917.23"
$ws.Range("E12").Value = "  +0.48%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.677.78"
$ws.Range("E13").Value = "  +0.24%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "4.13"
$ws.Range("E14").Value = "  +1.26%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.536"
$ws.Range("E15").Value = "  +2.08%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "66.16"
$ws.Range("E16").Value = "  +0.87%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "27.089.82"
$ws.Range("E17").Value = "  +0.51%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "237.94"
$ws.Range("E18").Value = "  +1.43%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "8.15"
$ws.Range("E19").Value = "  +1.34%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.0₃0743"
$ws.Range("E20").Value = "  +1.53%  "
$ws.Range("E21").Value = "  +0.06%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.49"
$ws.Range("E22").Value = "  +1.28%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "9.37"
$ws.Range("E23").Value = "  +2.40%  "
$ws.Range("E24").Value = "  -1.99%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "146.81"
$ws.Range("E25").Value = "  +0.72%  "
$ws.Range("E26").Value = "  +1.12%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "16.33"
$ws.Range("E27").Value = "  +2.34%  "
$ws.Range("E28").Value = "  +0.51%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.999"
$ws.Range("E29").Value = "  -0.12%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.0499"
$ws.Range("E30").Value = "  +0.31%  "
$ws.Range("E31").Value = "  +0.27%  "
$ws.Range("B32").Value = "Filecoin"
$ws.Range("C32").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.36"
$ws.Range("E32").Value = "  +1.08%  "
$ws.Range("B33").Value = "Maker"
$ws.Range("C33").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.556.65"
$ws.Range("E33").Value = "  +5.86%  "
$ws.Range("E34").Value = "  +1.82%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.70"
$ws.Range("E35").Value = "  +2.61%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.605"
$ws.Range("E36").Value = "  +4.67%  "
$ws.Range("B37").Value = "ARBITRUM"
$ws.Range("C37").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.938"
$ws.Range("E37").Value = "  +4.77%  "
$ws.Range("B38").Value = "HuobiToken"
$ws.Range("C38").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.39"
$ws.Range("E38").Value = "  -1.14%  "
$ws.Range("E39").Value = "  +2.01%  "
$ws.Range("E40").Value = "  +1.30%  "
$ws.Range("E41").Value = "  +0.05%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "68.22"
$ws.Range("E42").Value = "  +2.31%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "5.64"
$ws.Range("E43").Value = "  -2.49%  "
$ws.Range("E44").Value = "  -1.72%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.825.46"
$ws.Range("E45").Value = "  +0.83%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.783"
$ws.Range("E46").Value = "  +0.66%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "90.80"
$ws.Range("E47").Value = "  +0.32%  "
$ws.Range("B48").Value = "BabyDogeCoin"
$ws.Range("C48").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0₆0108"
$ws.Range("E48").Value = "  +1.99%  "
$ws.Range("B49").Value = "RenderToken"
$ws.Range("C49").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.57"
$ws.Range("E49").Value = "  +2.49%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.105"
$ws.Range("E50").Value = "  +3.61%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "8.02"
$ws.Range("E51").Value = "  +4.67%  "
